$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 278.56
$ws.Range("I15").Value = 278.56
$ws.Range("K15").Value = 835.6800000000001
$ws.Range("M15").Value = -666.6800000000001
$ws.Range("H107").Value = 1109.9333
$ws.Range("I107").Value = 1004.0833
$ws.Range("J107").Value = 1533.3334
$ws.Range("K107").Value = 1004.0833
$ws.Range("L107").Value = 1533.3334
$ws.Range("M107").Value = 915.9167
$ws.Range("N107").Value = -5373.3334
$ws.Range("H137").Value = 1929.3214
$ws.Range("I137").Value = 1748.84
$ws.Range("J137").Value = 3433.3333
$ws.Range("K137").Value = 5246.52
$ws.Range("L137").Value = 10299.9999
$ws.Range("M137").Value = -2696.52
$ws.Range("N137").Value = -15399.9999
$ws.Range("H138").Value = 14288428
$ws.Range("I138").Value = 47619740
$ws.Range("K138").Value = 142859220
$ws.Range("M138").Value = -142854080

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3731.897
$ws.Range("I32").Value = 2457.8
$ws.Range("J32").Value = 20113.143
$ws.Range("K32").Value = 2457.8
$ws.Range("L32").Value = 20113.143
$ws.Range("M32").Value = -2170.8
$ws.Range("N32").Value = -20687.143
$ws.Range("H61").Value = 644589.75
$ws.Range("I61").Value = 693715.9
$ws.Range("J61").Value = 5950
$ws.Range("K61").Value = 693715.9
$ws.Range("L61").Value = 5950
$ws.Range("M61").Value = -693503.9
$ws.Range("N61").Value = -6374
$ws.Range("H132").Value = 16412.941
$ws.Range("I132").Value = 1995.84
$ws.Range("J132").Value = 56460.445
$ws.Range("K132").Value = 5987.52
$ws.Range("L132").Value = 169381.335
$ws.Range("M132").Value = -3457.52
$ws.Range("N132").Value = -174441.335
$ws.Range("H136").Value = 644589.75
$ws.Range("I136").Value = 693715.9
$ws.Range("J136").Value = 5950
$ws.Range("K136").Value = 2081147.7
$ws.Range("L136").Value = 17850
$ws.Range("M136").Value = -2078597.7
$ws.Range("N136").Value = -22950

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1525.4117
$ws.Range("I20").Value = 1810.3
$ws.Range("K20").Value = 1810.3
$ws.Range("M20").Value = -1563.3
$ws.Range("H86").Value = 1474.75
$ws.Range("I86").Value = 1328.3572
$ws.Range("J86").Value = 1730.9375
$ws.Range("K86").Value = 1328.3572
$ws.Range("L86").Value = 1730.9375
$ws.Range("M86").Value = -205.3571999999999
$ws.Range("N86").Value = -3976.9375
$ws.Range("H89").Value = 1474.75
$ws.Range("I89").Value = 1328.3572
$ws.Range("J89").Value = 1730.9375
$ws.Range("K89").Value = 6641.786
$ws.Range("L89").Value = 8654.6875
$ws.Range("M89").Value = -1025.786
$ws.Range("N89").Value = -19886.6875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 783.1579
$ws.Range("I16").Value = 712.38464
$ws.Range("J16").Value = 936.5
$ws.Range("K16").Value = 712.38464
$ws.Range("L16").Value = 936.5
$ws.Range("M16").Value = -425.38464
$ws.Range("N16").Value = -1510.5
$ws.Range("H31").Value = 3856.1592
$ws.Range("I31").Value = 2449.1724
$ws.Range("K31").Value = 2449.1724
$ws.Range("M31").Value = -2154.1724
$ws.Range("H34").Value = 3856.1592
$ws.Range("I34").Value = 2449.1724
$ws.Range("K34").Value = 2449.1724
$ws.Range("M34").Value = -2247.1724
$ws.Range("H58").Value = 8563.046
$ws.Range("I58").Value = 755.2954999999999
$ws.Range("J58").Value = 24922.143
$ws.Range("K58").Value = 755.2954999999999
$ws.Range("L58").Value = 24922.143
$ws.Range("M58").Value = -552.2954999999999
$ws.Range("N58").Value = -25328.143
$ws.Range("H113").Value = 783.1579
$ws.Range("I113").Value = 712.38464
$ws.Range("J113").Value = 936.5
$ws.Range("K113").Value = 712.38464
$ws.Range("L113").Value = 936.5
$ws.Range("M113").Value = 1457.61536
$ws.Range("N113").Value = -5276.5
$ws.Range("H132").Value = 2193.8206
$ws.Range("I132").Value = 1732.3429
$ws.Range("J132").Value = 6231.75
$ws.Range("K132").Value = 5197.028700000001
$ws.Range("L132").Value = 18695.25
$ws.Range("M132").Value = -2667.028700000001
$ws.Range("N132").Value = -23755.25
$ws.Range("H136").Value = 8563.046
$ws.Range("I136").Value = 755.2954999999999
$ws.Range("J136").Value = 24922.143
$ws.Range("K136").Value = 2265.8865
$ws.Range("L136").Value = 74766.429
$ws.Range("M136").Value = 284.1135000000004
$ws.Range("N136").Value = -79866.429

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1574.2778
$ws.Range("I5").Value = 1249.125
$ws.Range("K5").Value = 3747.375
$ws.Range("M5").Value = -3635.375
$ws.Range("H46").Value = 1995
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1995
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5985
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -6167
$ws.Range("H62").Value = 5536.385
$ws.Range("J62").Value = 6710.778
$ws.Range("L62").Value = 20132.334
$ws.Range("N62").Value = -21504.334
$ws.Range("H65").Value = 5536.385
$ws.Range("J65").Value = 6710.778
$ws.Range("L65").Value = 60397.002
$ws.Range("N65").Value = -67261.00200000001
$ws.Range("H116").Value = 1623.4
$ws.Range("I116").Value = 1280
$ws.Range("J116").Value = 1737.8667
$ws.Range("K116").Value = 3840
$ws.Range("L116").Value = 5213.6001
$ws.Range("M116").Value = -398
$ws.Range("N116").Value = -12097.6001
$ws.Range("H122").Value = 421.42856
$ws.Range("I122").Value = 291.66666
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 2624.99994
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -174.9999399999997
$ws.Range("N122").Value = -15700
$ws.Range("H131").Value = 773.0306399999999
$ws.Range("J131").Value = 785.55914
$ws.Range("L131").Value = 2356.67742
$ws.Range("N131").Value = -12436.67742
$ws.Range("H135").Value = 1574.2778
$ws.Range("I135").Value = 1249.125
$ws.Range("K135").Value = 11242.125
$ws.Range("M135").Value = -8707.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 25555.555
$ws.Range("I70").Value = 50000
$ws.Range("J70").Value = 18571.428
$ws.Range("K70").Value = 50000
$ws.Range("L70").Value = 18571.428
$ws.Range("M70").Value = -49730
$ws.Range("N70").Value = -19111.428
$ws.Range("H73").Value = 25555.555
$ws.Range("I73").Value = 50000
$ws.Range("J73").Value = 18571.428
$ws.Range("K73").Value = 50000
$ws.Range("L73").Value = 18571.428
$ws.Range("M73").Value = -49064
$ws.Range("N73").Value = -20443.428
$ws.Range("H113").Value = 2649.4546
$ws.Range("I113").Value = 1677.7142
$ws.Range("K113").Value = 1677.7142
$ws.Range("M113").Value = 492.2858000000001
$ws.Range("H132").Value = 43955.617
$ws.Range("I132").Value = 6767.3
$ws.Range("J132").Value = 167916.67
$ws.Range("K132").Value = 20301.9
$ws.Range("L132").Value = 503750.01
$ws.Range("M132").Value = -17771.9
$ws.Range("N132").Value = -508810.01
